# Bugs corrigidos e adicionadas funcionalidades
#
# Adds a new data row (34) to the sheet with two brand new text values
# ("blabla" / "vai dar erro") which become new shared-string entries, and
# updates the current view/selection to reflect the newly added row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Activate()

# New row of data, following the same "tag replace" table layout used by
# the rest of the sheet (original tag | insert/carimbo type | page | dest tag).
$ws.Range("A34").Value = "blabla"
$ws.Range("B34").Value = "texto"
$ws.Range("C34").Value = 12
$ws.Range("D34").Value = "vai dar erro"

# Scroll the view down so the newly added row is visible, and move the
# active selection to the next empty row right below it.
$win = $excel.ActiveWindow
$win.ScrollRow = 16
$win.ScrollColumn = 1
$ws.Range("A35").Select()
